# Apply the "fund" (基金受益憑證) sheet edit:
#  - Row 1 changes from a stray data row into a proper header row
#    (name, owner, dealer, quantity, face_value, currency, total,
#     property_category, category, date, legislator_name, legislator_id,
#     source_file, index).
#  - Rows 2-5 gain the same trailing metadata columns (I:O) that the other
#    property sheets already carry: property_category="fund",
#    category="normal", date="2011-11-22", legislator_name="費鴻泰",
#    legislator_id=1365, source_file="tmp1afe1", index=<row's own index>.
#  - Row 3's quantity cell (E3) was stored as the text "21495.90"; it becomes
#    the number 21495.9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基金受益憑證")

# --- Fix header row (row 1) -------------------------------------------------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "owner"
$ws.Cells.Item(1, 4).Value = "dealer"
$ws.Cells.Item(1, 5).Value = "quantity"
$ws.Cells.Item(1, 6).Value = "face_value"
$ws.Cells.Item(1, 7).Value = "currency"
$ws.Cells.Item(1, 8).Value = "total"
$ws.Cells.Item(1, 9).Value = "property_category"
$ws.Cells.Item(1, 10).Value = "category"
$ws.Cells.Item(1, 11).Value = "date"
$ws.Cells.Item(1, 12).Value = "legislator_name"
$ws.Cells.Item(1, 13).Value = "legislator_id"
$ws.Cells.Item(1, 14).Value = "source_file"
$ws.Cells.Item(1, 15).Value = "index"

# New trailing header cells (I1:O1) - copy formatting from an existing
# header cell (B1) so they pick up the same bold/border/centered style.
$ws.Cells.Item(1, 2).Copy() | Out-Null
$headerRange = $ws.Range($ws.Cells.Item(1, 9), $ws.Cells.Item(1, 15))
$headerRange.PasteSpecial(-4122) | Out-Null

# --- Row 3: quantity stored as text -> numeric ------------------------------
$ws.Cells.Item(3, 5).Value = 21495.9

# --- Add trailing metadata columns (I:O) for the 4 data rows ---------------
$rows = @(2, 3, 4, 5)
$indexes = @(94, 95, 96, 97)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $idx = $indexes[$i]

    $ws.Cells.Item($r, 9).Value = "fund"
    $ws.Cells.Item($r, 10).Value = "normal"
    # Force the date column to stay literal text ("2011-11-22") instead of
    # being auto-converted to an Excel date serial number.
    $ws.Cells.Item($r, 11).NumberFormat = "@"
    $ws.Cells.Item($r, 11).Value = "2011-11-22"
    $ws.Cells.Item($r, 12).Value = "費鴻泰"
    $ws.Cells.Item($r, 13).Value = 1365
    $ws.Cells.Item($r, 14).Value = "tmp1afe1"
    $ws.Cells.Item($r, 15).Value = $idx
}

# Copy the data-row style (from an existing data cell, e.g. B2) onto the new
# cells so they match the rest of the row (same font/border/number format as
# the other data cells in that row).
$ws.Cells.Item(2, 2).Copy() | Out-Null
$dataRange = $ws.Range($ws.Cells.Item(2, 9), $ws.Cells.Item(5, 15))
$dataRange.PasteSpecial(-4122) | Out-Null
